# Actualización 10 de Mayo
# Update the statistics cells on "Estadisticos 1P", "Estadisticos 2P" and
# "Estadisticos Final" sheets (Blancos / Aprobados / Por_Apro / Promedio
# columns D, F, G, H) with the revised grading numbers.

$wb = $excel.ActiveWorkbook

# --- Estadisticos 1P ---------------------------------------------------
$ws1 = $wb.Worksheets.Item("Estadisticos 1P")
$ws1.Range("D2").Value = 6
$ws1.Range("F2").Value = 33
$ws1.Range("G2").Value = 84.62
$ws1.Range("H2").Value = 9.3

# --- Estadisticos 2P ---------------------------------------------------
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")

$ws2.Range("D2").Value = 13
$ws2.Range("F2").Value = 26
$ws2.Range("G2").Value = 66.67
$ws2.Range("H2").Value = 9.6

$ws2.Range("D3").Value = 8
$ws2.Range("E3").Value = 4
$ws2.Range("F3").Value = 27
$ws2.Range("G3").Value = 77.14
$ws2.Range("H3").Value = 9.4

$ws2.Range("D4").Value = 8
$ws2.Range("E4").Value = 4
$ws2.Range("F4").Value = 27
$ws2.Range("G4").Value = 77.14
$ws2.Range("H4").Value = 9.1

# --- Estadisticos Final -------------------------------------------------
$ws3 = $wb.Worksheets.Item("Estadisticos Final")

$ws3.Range("D2").Value = 6
$ws3.Range("F2").Value = 33
$ws3.Range("G2").Value = 84.62
$ws3.Range("H2").Value = 9.4

$ws3.Range("H4").Value = 9.5
